$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E width change (target OOXML width 10.453125; engine quantizes
# ColumnWidth to steps of 1/6, so 9.666666666666666 -> stored width 10.5,
# the closest representable value to the target)
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666

# Row 2 updates
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.023917979180256017
$ws.Range("D2").Value = 60.175010850415752
$ws.Range("E2").Value = 100.46847197090533
$ws.Range("F2").Value = 0.012710409572190263
$ws.Range("G2").Value = 60.160343854906323
$ws.Range("H2").Value = 100.4495788777019
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 7.754907353445911
$ws.Range("K2").Value = 38.879999999999995
$ws.Range("L2").Value = 0.19945749365858828

# Row 3 updates
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0.017717331897139133
$ws.Range("D3").Value = 62.191760722286176
$ws.Range("E3").Value = 96.70264579205886
$ws.Range("F3").Value = 0.0066459124260958838
$ws.Range("G3").Value = 62.179215652933962
$ws.Range("H3").Value = 96.690134793517515
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 7.2129768560020224
$ws.Range("K3").Value = 38.159999999999997
$ws.Range("L3").Value = 0.18901930964365887
